# Update NATMI LR-pair TPM stats (Fgf2-Cd44) with new TPM-derived values.
# Columns: E=Ligand-expressing cells, F=Ligand detection rate,
#          G=Ligand avg expr, H=Ligand total expr,
#          I=Ligand derived specificity (avg), J=Ligand derived specificity (total),
#          M=Receptor avg expr, N=Receptor total expr,
#          O=Receptor derived specificity (avg), P=Receptor derived specificity (total),
#          Q=Edge avg expr weight, R=Edge total expr weight,
#          S=Edge avg expr derived specificity, T=Edge total expr derived specificity

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$updates = @{
    2  = @{ E=2; F=0.6666666666666666; G=0.456404; H=1.369212; I=0.01914960767004715; J=0.01914960767004715;
            M=8.142376000000001; N=24.427128; O=0.1741313933276368; P=0.1741313933276368;
            Q=3.716212975904; R=33.44591678313601; S=0.003334547865262911; T=0.003334547865262911 }
    3  = @{ E=2; F=0.6666666666666666; G=0.456404; H=1.369212; I=0.01914960767004715; J=0.01914960767004715;
            O=0.5205382400466131; P=0.5205382400466131;
            Q=11.10903051511067; R=99.98127463599602; S=0.009968103074149469; T=0.009968103074149469 }
    4  = @{ E=2; F=0.6666666666666666; G=0.456404; H=1.369212; I=0.01914960767004715; J=0.01914960767004715;
            O=0.3053303666257501; P=0.3053303666257501;
            Q=6.516186706536001; R=58.645680358824; S=0.005846956730634774; T=0.005846956730634774 }
    5  = @{ I=0.8285024587002443; J=0.8285024587002443;
            M=8.142376000000001; N=24.427128; O=0.1741313933276368; P=0.1741313933276368;
            Q=160.7809225463174; R=1447.028302916856; S=0.1442682875088464; T=0.1442682875088464 }
    6  = @{ I=0.8285024587002443; J=0.8285024587002443;
            O=0.5205382400466131; P=0.5205382400466131;
            S=0.4312672117261169; T=0.4312672117261169 }
    7  = @{ I=0.8285024587002443; J=0.8285024587002443;
            O=0.3053303666257501; P=0.3053303666257501;
            S=0.252966959465281; T=0.252966959465281 }
    8  = @{ I=0.1523479336297086; J=0.1523479336297086;
            M=8.142376000000001; N=24.427128; O=0.1741313933276368; P=0.1741313933276368;
            Q=29.56495911362401; R=266.0846320226161; S=0.0265285579535275; T=0.0265285579535275 }
    9  = @{ I=0.1523479336297086; J=0.1523479336297086;
            O=0.5205382400466131; P=0.5205382400466131;
            S=0.07930292524634676; T=0.07930292524634676 }
    10 = @{ I=0.1523479336297086; J=0.1523479336297086;
            O=0.3053303666257501; P=0.3053303666257501;
            S=0.04651645042983438; T=0.04651645042983438 }
}

foreach ($row in $updates.Keys) {
    foreach ($col in $updates[$row].Keys) {
        $ws.Range("$col$row").Value = $updates[$row][$col]
    }
}
